$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5533
$ws1.Range("F9").Value = 6
$ws1.Range("F10").Value = 1079
$ws1.Range("F12").Value = 1541
$ws1.Range("F13").Value = 5080
$ws1.Range("F14").Value = 455
$ws1.Range("F16").Value = 209
$ws1.Range("F17").Value = 22
$ws1.Range("F18").Value = 8
$ws1.Range("F20").Value = 4362
$ws1.Range("F21").Value = 209
$ws1.Range("F22").Value = 1155
$ws1.Range("F24").Value = 61
$ws1.Range("F26").Value = 57
$ws1.Range("F27").Value = 166
$ws1.Range("F32").Value = 11
$ws1.Range("F37").Value = 43

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 208
$ws4.Range("F4").Value = 5533
$ws4.Range("F10").Value = 6
$ws4.Range("F11").Value = 1079
$ws4.Range("F13").Value = 1541
$ws4.Range("F14").Value = 5080
$ws4.Range("F15").Value = 455
$ws4.Range("F17").Value = 209
$ws4.Range("F18").Value = 22
$ws4.Range("F19").Value = 8
$ws4.Range("F21").Value = 4362
$ws4.Range("F22").Value = 209
$ws4.Range("F23").Value = 1155
$ws4.Range("F25").Value = 61
$ws4.Range("F27").Value = 57
$ws4.Range("F28").Value = 166
$ws4.Range("F33").Value = 11
$ws4.Range("F38").Value = 43
